$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-53 down to 29-54.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 3
$ws.Cells.Item(28, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44895
$ws.Cells.Item(28, 5).Value = 5
$ws.Cells.Item(28, 6).Value = 300000000
$ws.Cells.Item(28, 7).Value = "Espárragos"
$ws.Cells.Item(28, 8).Value = "Verde"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 1300
$ws.Cells.Item(28, 11).Value = 1500
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 13).Value = 1500
$ws.Cells.Item(28, 14).Value = "`$/kilo"
$ws.Cells.Item(28, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(28, 16).Value = 1500
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = "Hortaliza"
